# LOG.xlsx update: "python book, kicad plugin"
# - Added python book
# - installed kicad plugin that saves hours of time finding LCSC parts for
#   pcb assembly (kicad jlcpcb tools)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 28: fill in end time (D28) -> 16:00 ---------------------------
$ws.Range("C28").Copy()
$ws.Range("D28").PasteSpecial(-4122)   # xlPasteFormats (keep existing time style)
$excel.CutCopyMode = $false
$ws.Range("D28").Value = 0.66666666666666663

# --- Row 29: new log entry (date + start time) --------------------------
$ws.Range("B28").Copy()
$ws.Range("B29").PasteSpecial(-4122)   # xlPasteFormats (keep existing date style)
$excel.CutCopyMode = $false
$ws.Range("B29").Value = 44838

$ws.Range("C28").Copy()
$ws.Range("C29").PasteSpecial(-4122)   # xlPasteFormats (keep existing time style)
$excel.CutCopyMode = $false
$ws.Range("C29").Value = 0.42708333333333331

# --- Rows 82-87: clear out the leftover "C-number" text values ----------
$ws.Range("D82:D87").ClearContents()

# --- Update selection to reflect where the user was working -------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D29").Select()
